# Update CDD document version 1.4
# Fill in review decision / acceptance / comment columns on the review sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: "CDD file shall be "KeyPad CDD" not "CDD"" -> Accepted
$ws.Range("E2").Value = "Accepted"

# Row 3: "Init fucntion has no flow graph" -> Rejected, with a rejection comment
$ws.Range("E3").Value = "Rejected"
$ws.Range("G3").Value = "init function used only when there's global variables and in my application there's no need to use it "
$ws.Rows.Item(3).EntireRow.AutoFit()

# Row 4: "Sequence diagram..." -> Accepted
$ws.Range("E4").Value = "Accepted"

# Move the active selection to G4 as in the authored edit
$ws.Range("G4").Select()
